$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 7 (Ano 2025) with refreshed faturamento data
$ws.Range("B7").Value = 2957181.43
$ws.Range("C7").Value = -33.44295665906499
$ws.Range("D7").Value = 3003
$ws.Range("E7").Value = 3003
$ws.Range("F7").Value = 984.742400932401
$ws.Range("G7").Value = 4.966419334887839
